$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
   WHERE c.race = "NOT_REPORTED"
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@

$filesQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
        WHERE c.race = "NOT_REPORTED"
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@

$statQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
      WHERE c.race = "NOT_REPORTED"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

# Add new FilesTab row (row 3) label first
$ws.Range("A3").Value2 = "FilesTab"

# Update existing CasesTab row (row 2) with the new query text
$ws.Range("B2").Value2 = $casesQuery

# Fill in row 3 query cells
$ws.Range("B3").Value2 = $filesQuery

# Update the stat query cells (shared between rows 2 and 3)
$ws.Range("C2").Value2 = $statQuery
$ws.Range("C3").Value2 = $statQuery

$ws.Range("D3").Value2 = $ws.Range("D2").Value2
$ws.Range("E3").Value2 = $ws.Range("E2").Value2
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

$ws.Range("B2").RowHeight = 195
$ws.Range("B3").RowHeight = 409.5

# Update view settings
$ws.Range("B3").Select()
$excel.ActiveWindow.Zoom = 70

Write-Output "done"
